# "add aug 18 incher"
# Insert a new Steyr AUG A3 457mm (18") barrel row into the m4-barrels sheet,
# between the existing 508mm (row 29) and 407mm (row 30) entries, and tweak a
# handful of neighboring stat values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tweak existing rows 28 and 29 stats -----------------------------------
# Row 28: Steyr AUG HBAR 620mm
$ws.Range("E28").Value = -6
$ws.Range("H28").Value = 0.1

# Row 29: Steyr AUG A3 508mm
$ws.Range("C29").Value = -4
$ws.Range("E29").Value = -3
$ws.Range("F29").Value = -4
$ws.Range("H29").Value = -0.15

# --- Insert the new row for the 457mm (18") barrel -------------------------
# Shifts the old row 30 (407mm) and everything below it down by one.
$ws.Rows("30").Insert()

$ws.Range("A30").Value = "steyr_aug_457mm_barrel"
$ws.Range("B30").Value = "Steyr AUG A3 5.56x45 457mm"
$ws.Range("C30").Value = -2
$ws.Range("D30").Value = 0.69
$ws.Range("E30").Value = -2
$ws.Range("F30").Value = -2
$ws.Range("H30").Value = -0.08
$ws.Range("I30").Value = 0.16
$ws.Range("J30").Value = 275
$ws.Range("M30").Value = 750
$ws.Range("N30").Formula = "=C30-D30*20-E30*0.8-F30*0.6-H30*5+I30*10+J30/300"
$ws.Range("P30").Value = 0.1
$ws.Range("Q30").Value = 17.9921
$ws.Range("S30").Formula = "=ROUND(Q30*0.033+P30+R30, 2)"

# --- Tweak the row that used to be 31 (now 32) after the shift -------------
# Steyr AUG A3 350mm
$ws.Range("F32").Value = 0
$ws.Range("I32").Value = 0

# --- Minor view-state cosmetics (zoom level / selected cell) ---------------
$ws.Application.ActiveWindow.Zoom = 84
$ws.Range("D19").Select()
$ws.PageSetup.Orientation = 1
